$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H109").Value = 32943.332
$ws.Range("J109").Value = 32943.332
$ws.Range("L109").Value = 32943.332
$ws.Range("N109").Value = -35717.332
$ws.Range("H131").Value = 2165.45
$ws.Range("I131").Value = 1251.8182
$ws.Range("J131").Value = 3282.111
$ws.Range("K131").Value = 3755.4546
$ws.Range("L131").Value = 9846.332999999999
$ws.Range("M131").Value = 1284.5454
$ws.Range("N131").Value = -19926.333
$ws.Range("H132").Value = 13965.223
$ws.Range("I132").Value = 2091.6406
$ws.Range("J132").Value = 108953.875
$ws.Range("K132").Value = 6274.9218
$ws.Range("L132").Value = 326861.625
$ws.Range("M132").Value = -3744.9218
$ws.Range("N132").Value = -331921.625
$ws.Range("H135").Value = 7693395.5
$ws.Range("I135").Value = 720.5909
$ws.Range("K135").Value = 6485.3181
$ws.Range("M135").Value = -3950.3181
$ws.Range("H137").Value = 2622.56
$ws.Range("I137").Value = 904.7547
$ws.Range("J137").Value = 6760.909
$ws.Range("K137").Value = 2714.2641
$ws.Range("L137").Value = 20282.727
$ws.Range("M137").Value = -164.2640999999999
$ws.Range("N137").Value = -25382.727
$ws.Range("H138").Value = 1368.2538
$ws.Range("I138").Value = 1116.4916
$ws.Range("J138").Value = 3225
$ws.Range("K138").Value = 3349.4748
$ws.Range("L138").Value = 9675
$ws.Range("M138").Value = 1790.5252
$ws.Range("N138").Value = -19955
$ws.Range("H141").Value = 1401.1428
$ws.Range("I141").Value = 982.6727
$ws.Range("J141").Value = 4278.125
$ws.Range("K141").Value = 2948.0181
$ws.Range("L141").Value = 12834.375
$ws.Range("M141").Value = 2231.9819
$ws.Range("N141").Value = -23194.375

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9447.986999999999
$ws.Range("I32").Value = 8077.1406
$ws.Range("J32").Value = 21614.25
$ws.Range("K32").Value = 8077.1406
$ws.Range("L32").Value = 21614.25
$ws.Range("M32").Value = -7790.1406
$ws.Range("N32").Value = -22188.25
$ws.Range("H80").Value = 49086.285
$ws.Range("J80").Value = 49086.285
$ws.Range("L80").Value = 49086.285
$ws.Range("N80").Value = -51082.285
$ws.Range("H83").Value = 49086.285
$ws.Range("J83").Value = 49086.285
$ws.Range("L83").Value = 147258.855
$ws.Range("N83").Value = -157242.855
$ws.Range("H102").Value = 11326.087
$ws.Range("I102").Value = 1856.6666
$ws.Range("J102").Value = 21656.363
$ws.Range("K102").Value = 1856.6666
$ws.Range("L102").Value = 21656.363
$ws.Range("M102").Value = -234.6666
$ws.Range("N102").Value = -24900.363
$ws.Range("H103").Value = 39362
$ws.Range("J103").Value = 39362
$ws.Range("L103").Value = 39362
$ws.Range("N103").Value = -41706
$ws.Range("H110").Value = 1635.9231
$ws.Range("I110").Value = 1596.2858
$ws.Range("J110").Value = 1802.4
$ws.Range("K110").Value = 1596.2858
$ws.Range("L110").Value = 1802.4
$ws.Range("M110").Value = 448.7141999999999
$ws.Range("N110").Value = -5892.4
$ws.Range("H122").Value = 2092.3447
$ws.Range("I122").Value = 2322.762
$ws.Range("K122").Value = 6968.286
$ws.Range("M122").Value = -4518.286
$ws.Range("H132").Value = 6945744
$ws.Range("I132").Value = 8772792
$ws.Range("J132").Value = 2963.0667
$ws.Range("K132").Value = 26318376
$ws.Range("L132").Value = 8889.2001
$ws.Range("M132").Value = -26315846
$ws.Range("N132").Value = -13949.2001

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H116").Value = 43380.5
$ws.Range("J116").Value = 43380.5
$ws.Range("L116").Value = 43380.5
$ws.Range("N116").Value = -52558.5
$ws.Range("H117").Value = 44257.832
$ws.Range("J117").Value = 44257.832
$ws.Range("L117").Value = 44257.832
$ws.Range("N117").Value = -53435.832
$ws.Range("H124").Value = 48984
$ws.Range("J124").Value = 48984
$ws.Range("L124").Value = 48984
$ws.Range("N124").Value = -58804
$ws.Range("H126").Value = 43606
$ws.Range("J126").Value = 43606
$ws.Range("L126").Value = 43606
$ws.Range("N126").Value = -53486
$ws.Range("H130").Value = 45081.8
$ws.Range("J130").Value = 45081.8
$ws.Range("L130").Value = 45081.8
$ws.Range("N130").Value = -55121.8
$ws.Range("H134").Value = 178421.12
$ws.Range("I134").Value = 1513.7916
$ws.Range("K134").Value = 4541.3748
$ws.Range("M134").Value = -2006.3748

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1270.6222
$ws.Range("I58").Value = 901.38464
$ws.Range("J58").Value = 1775.8948
$ws.Range("K58").Value = 901.38464
$ws.Range("L58").Value = 1775.8948
$ws.Range("M58").Value = -698.38464
$ws.Range("N58").Value = -2181.8948
$ws.Range("H111").Value = 49247
$ws.Range("J111").Value = 49247
$ws.Range("L111").Value = 49247
$ws.Range("N111").Value = -57427
$ws.Range("H118").Value = 49267
$ws.Range("J118").Value = 49267
$ws.Range("L118").Value = 49267
$ws.Range("N118").Value = -52581
$ws.Range("H122").Value = 121979.9
$ws.Range("I122").Value = 201516.5
$ws.Range("J122").Value = 2675
$ws.Range("K122").Value = 604549.5
$ws.Range("L122").Value = 8025
$ws.Range("M122").Value = -602099.5
$ws.Range("N122").Value = -12925
$ws.Range("H132").Value = 15648.876
$ws.Range("I132").Value = 1045.5975
$ws.Range("J132").Value = 95480.13
$ws.Range("K132").Value = 3136.7925
$ws.Range("L132").Value = 286440.39
$ws.Range("M132").Value = -606.7925000000005
$ws.Range("N132").Value = -291500.39
$ws.Range("H134").Value = 213228.11
$ws.Range("I134").Value = 795.3
$ws.Range("J134").Value = 877080.6
$ws.Range("K134").Value = 2385.9
$ws.Range("L134").Value = 2631241.8
$ws.Range("M134").Value = 149.1000000000004
$ws.Range("N134").Value = -2636311.8
$ws.Range("H136").Value = 1270.6222
$ws.Range("I136").Value = 901.38464
$ws.Range("J136").Value = 1775.8948
$ws.Range("K136").Value = 2704.15392
$ws.Range("L136").Value = 5327.6844
$ws.Range("M136").Value = -154.1539199999997
$ws.Range("N136").Value = -10427.6844

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 45474.855
$ws.Range("J130").Value = 45474.855
$ws.Range("L130").Value = 45474.855
$ws.Range("N130").Value = -55514.855

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 16500
$ws.Range("J50").Value = 16500
$ws.Range("L50").Value = 16500
$ws.Range("N50").Value = -17774
$ws.Range("H110").Value = 25500
$ws.Range("J110").Value = 25500
$ws.Range("L110").Value = 25500
$ws.Range("N110").Value = -33680
$ws.Range("H111").Value = 43756
$ws.Range("J111").Value = 43756
$ws.Range("L111").Value = 43756
$ws.Range("N111").Value = -51936

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1309.6279
$ws.Range("I132").Value = 669.125
$ws.Range("J132").Value = 2118.6843
$ws.Range("K132").Value = 2007.375
$ws.Range("L132").Value = 6356.0529
$ws.Range("M132").Value = 522.625
$ws.Range("N132").Value = -11416.0529
$ws.Range("H136").Value = 14256.311
$ws.Range("I136").Value = 17126.1
$ws.Range("J136").Value = 1957.2142
$ws.Range("K136").Value = 51378.3
$ws.Range("L136").Value = 5871.642599999999
$ws.Range("M136").Value = -48828.3
$ws.Range("N136").Value = -10971.6426
